$d = $word.ActiveDocument

# --- Step 1: Insert a new title paragraph before paragraph 1 ---
$p1 = $d.Paragraphs.Item(1)
$titlePara = $p1.Range.InsertParagraphBefore()

# The new (now) paragraph 1 is empty; fill it in with three runs.
$titleRange = $d.Paragraphs.Item(1).Range
$titleRange.InsertBefore("3 Day in situ Hi-C")

$titleStart = $d.Paragraphs.Item(1).Range.Start
$italicRange = $d.Range($titleStart + 6, $titleStart + 13)
$italicRange.Font.Italic = $true

Write-Host "After step1:"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "Para $i : [$($p.Range.Text)]"
}

# --- Step 2: rewrite old paragraph 1 (now paragraph 2) ---
# "I.a.2. In situ Hi-C libraries..." -> " " + italic("i") + italic("n situ") + " Hi-C libraries..."
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start
$target = $d.Range($p2Start, $p2Start + 8)   # "I.a.2. I"
Write-Host "Target to replace: [$($target.Text)]"
$target.Text = " i"

$p2Start2 = $d.Paragraphs.Item(2).Range.Start
$italicRange2 = $d.Range($p2Start2 + 1, $p2Start2 + 8)  # "in situ"
Write-Host "Italic target: [$($italicRange2.Text)]"
$italicRange2.Font.Italic = $true

Write-Host "After step2:"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "Para $i : [$($p.Range.Text)]"
}

# --- Step 3: rewrite paragraph 3 ("whether ...") to drop the proofErr markers ---
# Insert a clean paragraph after paragraph 2, move the text there, delete the old paragraph 3.
$p2b = $d.Paragraphs.Item(2)
$p2b.Range.InsertParagraphAfter() | Out-Null
$newP3 = $d.Paragraphs.Item(3).Range
$newP3.InsertBefore("whether shorter incubation times are used for the restriction step (2 hours; see step 12), the fill-in step (45 minutes; see")

# Delete the old paragraph (now paragraph 4, holding the proofErr markup)
$oldP3 = $d.Paragraphs.Item(4)
$oldP3full = $d.Range($oldP3.Range.Start, $oldP3.Range.End)
$oldP3full.Delete()

Write-Host "After step3:"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "Para $i : [$($p.Range.Text)]"
}
